$d = $word.ActiveDocument

function Replace-Exact($search, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $r.Text = $newText
    }
    return $found
}

# 1. "модную раскладку в карман" -> "модную раскладушку в карман"
Replace-Exact "раскладку в карман" "раскладушку в карман"

# 2. "Итак. С кем я имею честь?" -> "Итак, с кем я имею честь?"
Replace-Exact "Итак. С кем я имею честь?" "Итак, с кем я имею честь?"

# 3. Remove the stray _GoBack bookmark left over from the previous save
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 4. En dash after "руку." corrected to a plain hyphen
Replace-Exact "руку. – Очень любезно" "руку. - Очень любезно"

# 5. "говориться" -> "говорится" (spelling fix)
Replace-Exact "как говориться, и коту приятно" "как говорится, и коту приятно"

# 6. "сжал" -> "сжала" (gender agreement fix)
Replace-Exact "Я слегка сжал его руку" "Я слегка сжала его руку"

# 7. "клетку ... Вируса" -> "клетки ... вируса"
Replace-Exact "клетку для адаптации первого поколения Вируса гоблинов" "клетки для адаптации первого поколения вируса гоблинов"

# 8. Drop the trailing space after "Коннорса!" (and tidy the surrounding dash spacing)
Replace-Exact "Коннор! Отто чуть не закричал. - Это эксперимент Коннорса! " "Коннор! Отто чуть не закричал. - Это эксперимент Коннорса!"

# 9. Double space collapsed to a single space
Replace-Exact "пальцы  разжались" "пальцы разжались"

# 10. Drop trailing space at the end of the paragraph
Replace-Exact "Я знаю, чем занимается доктор Коннорс. " "Я знаю, чем занимается доктор Коннорс."

# 11. Remove the extra space between the closing quote and the ellipsis:
#     "Имобилизиен" ...  ->  "Имобилизиен"...
#     (quotes are ASCII straight quotes in the source; go through Range.Text
#      rather than Find.Execute's replacement argument so AutoCorrect's
#      straight-quotes-to-curly-quotes substitution is not triggered)
$r = $d.Content
$found = $r.Find.Execute([char]34 + "Имобилизиен" + [char]34 + " ...", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) {
    $r.Text = [char]34 + "Имобилизиен" + [char]34 + "..."
}

# 12. "гумаоида" -> "гуманоида" (spelling fix)
Replace-Exact "огромного гумаоида" "огромного гуманоида"
